$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.710097312927246
$ws.Range("B1").Value = 2.965962886810303
$ws.Range("C1").Value = 3.166210174560547
$ws.Range("D1").Value = 3.567677736282349
$ws.Range("E1").Value = 3.900121927261353
